$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(45995, 10792.687886729, 9964.17677749314, 16704.26, 6596.50352088445, -5.98248756759995)
    3 = @(45996, 10689.8613252878, 9325.15672308304, 11232.26, 6655.52991538388, 197.851109936121)
    4 = @(45997, 4053.37664031883, 6667.89865830184, 11232.26, 6367.29072722427, 75.1220577302548)
    5 = @(45998, 3857.76206559275, 6623.16198470777, 11232.26, 6304.69655355081, 70.6499390941074)
    6 = @(45999, 3828.41217004917, 6519.74269935672, 11232.26, 6338.63920467329, 67.755079334584)
    7 = @(46000, 9706.9828974431, 9708.3819725428, 11232.26, 7769.38200686065, 260.229332475144)
    8 = @(46001, 9706.9828974431, 9574.20127201797, 11232.26, 7769.38200686065, 254.638469953276)
    9 = @(46002, 9706.9828974431, 9576.82888712837, 11232.26, 7769.38200686065, 254.747953916209)
    10 = @(46003, 9706.9828974431, 8707.29935014787, 11232.26, 7769.38200686065, 218.517556542022)
    11 = @(46004, 8488.30621711765, 8455.5235119604, 11232.26, 7366.5902533736, 191.243906888917)
    12 = @(46005, 8386.73995792576, 8371.91465139008, 11232.26, 7357.47177409387, 187.380267728498)
    13 = @(46006, 9746.03784374469, 8926.05465055304, 11232.26, 8042.17402277399, 238.99869472196)
    14 = @(46007, 9746.03784374469, 9249.06507406336, 11232.26, 8042.17402277399, 252.457462368223)
    15 = @(46008, 9746.03784374469, 9164.5879239753, 11232.26, 8042.17402277399, 248.937581114554)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]
    $ws.Cells.Item($r, 5).Value = $vals[4]
    $ws.Cells.Item($r, 6).Value = $vals[5]
}

Write-Output "Done updating need_to_buy rows 2-15"
